$d = $word.ActiveDocument

# Locate the paragraph that ends with "...Windows dependencies whatsoever."
$findRange = $d.Content
$findRange.Find.Execute("avoiding any kind of Windows dependencies whatsoever.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorPara = $findRange.Paragraphs(1)
$anchorIndex = $anchorPara.Index

# Insert a brand-new empty paragraph right after the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter() | Out-Null

# The freshly inserted paragraph is now the one immediately after the anchor.
$newPara = $d.Paragraphs($anchorIndex + 1)

$newParagraphXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">You can read how-to articles in the growing library </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">at </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>Docs/Usage/</w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>HowTo</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>. If you don't find what you need right away, either check back soon or create an I</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">ssue </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">about the subject </w:t>
            </w:r>
            <w:r>
              <w:t>and I will make sure it gets explained.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$newPara.Range.InsertXML($newParagraphXml) | Out-Null
